$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.110.86"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "2.652.16"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'519.33"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").Value = "'147.15"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "2.669.14"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "'0.128"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "3.114.14"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "59.121.12"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "2.660.18"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "'352.01"
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'10.44"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").Value = "'6.23"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'61.89"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "0.0₃0811"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "'7.13"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'6.35"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").Value = "'19.00"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("D34").Value = "'149.76"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "'0.967"
$ws.Range("E35").Value = "  -5.88%  "
$ws.Range("D36").Value = "'4.06"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("D37").Value = "'1.19"
$ws.Range("E37").Value = "  +5.05%  "
$ws.Range("D38").Value = "'0.859"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "'36.66"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").Value = "'3.67"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Value = "'282.70"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.0991"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "'19.96"
$ws.Range("E45").Value = "  +3.25%  "
$ws.Range("D46").Value = "'0.608"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "2.109.39"
$ws.Range("E47").Value = "  +7.85%  "
$ws.Range("D48").Value = "'0.0531"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.81"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0232"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'10.31"
$ws.Range("E51").Value = "  +0.75%  "
